$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 6687
$ws.Range("I111").Value = 4231.6
$ws.Range("J111").Value = 11597.8
$ws.Range("K111").Value = 12694.8
$ws.Range("L111").Value = 34793.39999999999
$ws.Range("M111").Value = -9627.800000000001
$ws.Range("N111").Value = -40927.39999999999

$ws.Range("H132").Value = 5185.7144
$ws.Range("I132").Value = 1726.7826
$ws.Range("K132").Value = 5180.3478
$ws.Range("M132").Value = -2650.3478

$ws.Range("H137").Value = 24527778
$ws.Range("I137").Value = 1115006.5
$ws.Range("K137").Value = 3345019.5
$ws.Range("M137").Value = -3342469.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15039.237
$ws.Range("I32").Value = 14024.444
$ws.Range("J32").Value = 25999
$ws.Range("K32").Value = 14024.444
$ws.Range("L32").Value = 25999
$ws.Range("M32").Value = -13737.444
$ws.Range("N32").Value = -26573

$ws.Range("H45").Value = 3960
$ws.Range("J45").Value = 4450
$ws.Range("L45").Value = 4450
$ws.Range("N45").Value = -5204

$ws.Range("H61").Value = 4324.5835
$ws.Range("I61").Value = 4481.364
$ws.Range("K61").Value = 4481.364
$ws.Range("M61").Value = -4269.364

$ws.Range("H74").Value = 35717108
$ws.Range("I74").Value = 41669540
$ws.Range("J74").Value = 2500
$ws.Range("K74").Value = 41669540
$ws.Range("L74").Value = 2500
$ws.Range("M74").Value = -41668666
$ws.Range("N74").Value = -4248

$ws.Range("H77").Value = 35717108
$ws.Range("I77").Value = 41669540
$ws.Range("J77").Value = 2500
$ws.Range("K77").Value = 208347700
$ws.Range("L77").Value = 12500
$ws.Range("M77").Value = -208343332
$ws.Range("N77").Value = -21236

$ws.Range("H97").Value = 1820
$ws.Range("I97").Value = 1774.0667
$ws.Range("K97").Value = 1774.0667
$ws.Range("M97").Value = -1278.0667

$ws.Range("H110").Value = 1785.7142
$ws.Range("I110").Value = 1785.7142
$ws.Range("K110").Value = 1785.7142
$ws.Range("M110").Value = 259.2858000000001

$ws.Range("H132").Value = 1895.3334
$ws.Range("I132").Value = 1803.3055
$ws.Range("K132").Value = 5409.916499999999
$ws.Range("M132").Value = -2879.916499999999

$ws.Range("H136").Value = 4324.5835
$ws.Range("I136").Value = 4481.364
$ws.Range("K136").Value = 13444.092
$ws.Range("M136").Value = -10894.092

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2997.5
$ws.Range("I20").Value = 2997.5
$ws.Range("K20").Value = 2997.5
$ws.Range("M20").Value = -2750.5

$ws.Range("H86").Value = 2442.75
$ws.Range("I86").Value = 2159.6667
$ws.Range("J86").Value = 2725.8333
$ws.Range("K86").Value = 2159.6667
$ws.Range("L86").Value = 2725.8333
$ws.Range("M86").Value = -1036.6667
$ws.Range("N86").Value = -4971.8333

$ws.Range("H89").Value = 2442.75
$ws.Range("I89").Value = 2159.6667
$ws.Range("J89").Value = 2725.8333
$ws.Range("K89").Value = 10798.3335
$ws.Range("L89").Value = 13629.1665
$ws.Range("M89").Value = -5182.333500000001
$ws.Range("N89").Value = -24861.1665

$ws.Range("H94").Value = 944.0625
$ws.Range("I94").Value = 842.9091
$ws.Range("K94").Value = 842.9091
$ws.Range("M94").Value = -391.9091

$ws.Range("H107").Value = 2930.1333
$ws.Range("I107").Value = 2988.5
$ws.Range("J107").Value = 2813.4
$ws.Range("K107").Value = 2988.5
$ws.Range("L107").Value = 2813.4
$ws.Range("M107").Value = -1068.5
$ws.Range("N107").Value = -6653.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 600
$ws.Range("I86").Value = 600
$ws.Range("K86").Value = 600
$ws.Range("M86").Value = 523

$ws.Range("H89").Value = 600
$ws.Range("I89").Value = 600
$ws.Range("K89").Value = 3000
$ws.Range("M89").Value = 2616

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H113").Value = 754.8
$ws.Range("J113").Value = 727.6667
$ws.Range("L113").Value = 2183.0001
$ws.Range("N113").Value = -6523.0001

$ws.Range("H128").Value = 128499.5
$ws.Range("I128").Value = 128499.5
$ws.Range("K128").Value = 385498.5
$ws.Range("M128").Value = -380518.5

$ws.Range("H131").Value = 7733740
$ws.Range("I131").Value = 93145.45
$ws.Range("J131").Value = 9643889
$ws.Range("K131").Value = 279436.35
$ws.Range("L131").Value = 28931667
$ws.Range("M131").Value = -274396.35
$ws.Range("N131").Value = -28941747

$ws.Range("H140").Value = 2404.1428
$ws.Range("I140").Value = 2161.7693
$ws.Range("K140").Value = 6485.3079
$ws.Range("M140").Value = -1305.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8222.68
$ws.Range("I70").Value = 7948.3335
$ws.Range("K70").Value = 7948.3335
$ws.Range("M70").Value = -7678.3335

$ws.Range("H73").Value = 8222.68
$ws.Range("I73").Value = 7948.3335
$ws.Range("K73").Value = 7948.3335
$ws.Range("M73").Value = -7012.3335

$ws.Range("H93").Value = 44333.332
$ws.Range("J93").Value = 44333.332
$ws.Range("L93").Value = 44333.332
$ws.Range("N93").Value = -48077.332

$ws.Range("H113").Value = 900
$ws.Range("J113").Value = 600
$ws.Range("L113").Value = 600
$ws.Range("N113").Value = -4940

$ws.Range("H126").Value = 2265.04
$ws.Range("J126").Value = 3066.3333
$ws.Range("L126").Value = 9198.999899999999
$ws.Range("N126").Value = -14138.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1749.5
$ws.Range("J2").Value = 1749.5
$ws.Range("L2").Value = 1749.5
$ws.Range("N2").Value = -1973.5

$ws.Range("H16").Value = 3625.25
$ws.Range("I16").Value = 3625.25
$ws.Range("K16").Value = 3625.25
$ws.Range("M16").Value = -3455.25

$ws.Range("H122").Value = 7202.486
$ws.Range("I122").Value = 4092.5217
$ws.Range("K122").Value = 12277.5651
$ws.Range("M122").Value = -9827.5651

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 448.0909
$ws.Range("I100").Value = 473
$ws.Range("J100").Value = 199
$ws.Range("K100").Value = 946
$ws.Range("L100").Value = 398
$ws.Range("M100").Value = -405
$ws.Range("N100").Value = -1480

$ws.Range("H113").Value = 330.72
$ws.Range("I113").Value = 381.6
$ws.Range("K113").Value = 1144.8
$ws.Range("M113").Value = 1025.2
